$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(2, 6).Value = 2889
$ws1.Cells.Item(3, 6).Value = 21217
$ws1.Cells.Item(4, 6).Value = 102
$ws1.Cells.Item(5, 6).Value = 2937
$ws1.Cells.Item(7, 6).Value = 619
$ws1.Cells.Item(8, 6).Value = 518
$ws1.Cells.Item(9, 6).Value = 767
$ws1.Cells.Item(10, 6).Value = 280
$ws1.Cells.Item(13, 6).Value = 119
$ws1.Cells.Item(14, 6).Value = 522
$ws1.Cells.Item(15, 6).Value = 183
$ws1.Cells.Item(16, 6).Value = 269
$ws1.Cells.Item(17, 6).Value = 19
$ws1.Cells.Item(18, 6).Value = 423
$ws1.Cells.Item(19, 6).Value = 66
$ws1.Cells.Item(22, 6).Value = 36
$ws1.Cells.Item(23, 6).Value = 127

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item(2)
$ws2.Cells.Item(3, 6).Value = 33
$ws2.Cells.Item(5, 6).Value = 340
$ws2.Cells.Item(8, 6).Value = 18
$ws2.Cells.Item(14, 6).Value = 155

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 6132
$ws3.Cells.Item(3, 6).Value = 701
$ws3.Cells.Item(4, 6).Value = 698
$ws3.Cells.Item(5, 6).Value = 1621
$ws3.Cells.Item(6, 6).Value = 56

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(2, 6).Value = 6132
$ws4.Cells.Item(3, 6).Value = 701
$ws4.Cells.Item(4, 6).Value = 698
$ws4.Cells.Item(5, 6).Value = 1621
$ws4.Cells.Item(6, 6).Value = 2889
$ws4.Cells.Item(7, 6).Value = 21217
$ws4.Cells.Item(9, 6).Value = 33
$ws4.Cells.Item(10, 6).Value = 102
$ws4.Cells.Item(12, 6).Value = 340
$ws4.Cells.Item(13, 6).Value = 2937
$ws4.Cells.Item(16, 6).Value = 56
$ws4.Cells.Item(17, 6).Value = 619
$ws4.Cells.Item(18, 6).Value = 518
$ws4.Cells.Item(19, 6).Value = 767
$ws4.Cells.Item(20, 6).Value = 280
$ws4.Cells.Item(24, 6).Value = 18
$ws4.Cells.Item(26, 6).Value = 119
$ws4.Cells.Item(29, 6).Value = 522
$ws4.Cells.Item(31, 6).Value = 183
$ws4.Cells.Item(33, 6).Value = 269
$ws4.Cells.Item(34, 6).Value = 155
$ws4.Cells.Item(35, 6).Value = 155
$ws4.Cells.Item(36, 6).Value = 19
$ws4.Cells.Item(37, 6).Value = 423
$ws4.Cells.Item(39, 6).Value = 66
$ws4.Cells.Item(44, 6).Value = 36
$ws4.Cells.Item(50, 6).Value = 127
